$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "nom" column (DR), shifting
# "nom" (DR -> DS) and "url_produit" (DS -> DT) one column to the right.
$ws.Columns("DR").Insert()

# New header cell for the freshly inserted price-history column.
$ws.Range("DR1").Value = "2026-02-02 08:29:42"

# For every data row, carry the last known price (column DQ, the most
# recent existing snapshot) forward into the newly inserted column DR.
# Rows whose price history is blank (already discontinued products) stay
# blank in the new column too.
for ($r = 2; $r -le 206; $r++) {
    $lastPrice = $ws.Cells.Item($r, 121).Value2
    if ($lastPrice -ne "") {
        $ws.Cells.Item($r, 122).Value = $lastPrice
    }
}
